# Auto-generated edit script: updates market-price derived columns (H-N)
# across multiple sheets to match the scheduled runner refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 1771.0769
$ws.Range("I48").Value = 1093.091
$ws.Range("K48").Value = 3279.273
$ws.Range("M48").Value = -2987.273
$ws.Range("H56").Value = 1771.0769
$ws.Range("I56").Value = 1093.091
$ws.Range("K56").Value = 3279.273
$ws.Range("M56").Value = -2745.273
$ws.Range("H98").Value = 1814.7097
$ws.Range("J98").Value = 3778.3333
$ws.Range("L98").Value = 3778.3333
$ws.Range("N98").Value = -6774.3333
$ws.Range("H100").Value = 2838.5789
$ws.Range("I100").Value = 1944.5385
$ws.Range("K100").Value = 1944.5385
$ws.Range("M100").Value = -1403.5385
$ws.Range("H106").Value = 15107.556
$ws.Range("I106").Value = 2997.25
$ws.Range("K106").Value = 2997.25
$ws.Range("M106").Value = -2366.25
$ws.Range("H113").Value = 5249.875
$ws.Range("J113").Value = 6199.8
$ws.Range("L113").Value = 6199.8
$ws.Range("N113").Value = -12707.8
$ws.Range("H122").Value = 1814.7097
$ws.Range("J122").Value = 3778.3333
$ws.Range("L122").Value = 11334.9999
$ws.Range("N122").Value = -16234.9999
$ws.Range("H132").Value = 1528.96
$ws.Range("I132").Value = 1331.4783
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 3994.4349
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -1464.4349
$ws.Range("N132").Value = -16460
$ws.Range("H134").Value = 134874.5
$ws.Range("I134").Value = 149749
$ws.Range("J134").Value = 120000
$ws.Range("K134").Value = 149749
$ws.Range("L134").Value = 120000
$ws.Range("M134").Value = -144679
$ws.Range("N134").Value = -130140
$ws.Range("H136").Value = 131875.67
$ws.Range("J136").Value = 131875.67
$ws.Range("L136").Value = 131875.67
$ws.Range("N136").Value = -142075.67
$ws.Range("H137").Value = 1928.55
$ws.Range("I137").Value = 1739.7646
$ws.Range("J137").Value = 2998.3333
$ws.Range("K137").Value = 5219.293799999999
$ws.Range("L137").Value = 8994.999899999999
$ws.Range("M137").Value = -2669.293799999999
$ws.Range("N137").Value = -14094.9999
$ws.Range("H138").Value = 2730.68
$ws.Range("I138").Value = 2498.4211
$ws.Range("J138").Value = 3466.1667
$ws.Range("K138").Value = 7495.263300000001
$ws.Range("L138").Value = 10398.5001
$ws.Range("M138").Value = -2355.263300000001
$ws.Range("N138").Value = -20678.5001
$ws.Range("H139").Value = 111189.75
$ws.Range("J139").Value = 111189.75
$ws.Range("L139").Value = 111189.75
$ws.Range("N139").Value = -121469.75
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18477.156
$ws.Range("I32").Value = 4252.203
$ws.Range("K32").Value = 4252.203
$ws.Range("M32").Value = -3965.203
$ws.Range("H41").Value = 2832.889
$ws.Range("I41").Value = 2832.889
$ws.Range("K41").Value = 2832.889
$ws.Range("M41").Value = -2418.889
$ws.Range("H102").Value = 2677.7856
$ws.Range("I102").Value = 1622.875
$ws.Range("K102").Value = 1622.875
$ws.Range("M102").Value = -0.875
$ws.Range("H132").Value = 2134.6956
$ws.Range("I132").Value = 2134.6956
$ws.Range("K132").Value = 6404.0868
$ws.Range("M132").Value = -3874.0868
$ws.Range("H133").Value = 96614.28999999999
$ws.Range("J133").Value = 96614.28999999999
$ws.Range("L133").Value = 96614.28999999999
$ws.Range("N133").Value = -101674.29
$ws.Range("H137").Value = 88980.336
$ws.Range("J137").Value = 88980.336
$ws.Range("L137").Value = 88980.336
$ws.Range("N137").Value = -99180.336
$ws.Range("H138").Value = 118974
$ws.Range("J138").Value = 118974
$ws.Range("L138").Value = 118974
$ws.Range("N138").Value = -129254
$ws.Range("H139").Value = 97199.2
$ws.Range("J139").Value = 97199.2
$ws.Range("L139").Value = 97199.2
$ws.Range("N139").Value = -107479.2

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8391.9
$ws.Range("I20").Value = 7854.846
$ws.Range("K20").Value = 7854.846
$ws.Range("M20").Value = -7607.846
$ws.Range("H68").Value = 25295
$ws.Range("J68").Value = 25295
$ws.Range("L68").Value = 25295
$ws.Range("N68").Value = -26917
$ws.Range("H71").Value = 25295
$ws.Range("J71").Value = 25295
$ws.Range("L71").Value = 75885
$ws.Range("N71").Value = -83997
$ws.Range("H99").Value = 3565.5925
$ws.Range("I99").Value = 3591
$ws.Range("K99").Value = 3591
$ws.Range("M99").Value = -2093

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1206.0244
$ws.Range("I107").Value = 1086.375
$ws.Range("J107").Value = 1374.9412
$ws.Range("K107").Value = 1086.375
$ws.Range("L107").Value = 1374.9412
$ws.Range("M107").Value = 833.625
$ws.Range("N107").Value = -5214.9412
$ws.Range("H125").Value = 49163
$ws.Range("J125").Value = 49163
$ws.Range("L125").Value = 49163
$ws.Range("N125").Value = -54083
$ws.Range("H132").Value = 3323.9167
$ws.Range("I132").Value = 3431.889
$ws.Range("K132").Value = 10295.667
$ws.Range("M132").Value = -7765.667000000001
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 287724.72
$ws.Range("J141").Value = 287724.72
$ws.Range("L141").Value = 287724.72
$ws.Range("N141").Value = -298084.72

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 5001499
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H60").Value = 1049.5
$ws.Range("I60").Value = 1049.5
$ws.Range("K60").Value = 3148.5
$ws.Range("M60").Value = -2897.5
$ws.Range("H107").Value = 818
$ws.Range("I107").Value = 680.5
$ws.Range("J107").Value = 852.375
$ws.Range("K107").Value = 2041.5
$ws.Range("L107").Value = 2557.125
$ws.Range("M107").Value = -121.5
$ws.Range("N107").Value = -6397.125
$ws.Range("H121").Value = 18520664
$ws.Range("J121").Value = 2675.8572
$ws.Range("L121").Value = 8027.571599999999
$ws.Range("N121").Value = -10647.5716
$ws.Range("H131").Value = 12397.741
$ws.Range("I131").Value = 836.2857
$ws.Range("J131").Value = 24848.54
$ws.Range("K131").Value = 2508.8571
$ws.Range("L131").Value = 74545.62
$ws.Range("M131").Value = 2531.1429
$ws.Range("N131").Value = -84625.62
$ws.Range("H132").Value = 1739.7826
$ws.Range("J132").Value = 2067.6428
$ws.Range("L132").Value = 18608.7852
$ws.Range("N132").Value = -23668.7852

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3268.5122
$ws.Range("I132").Value = 2070.4666
$ws.Range("J132").Value = 6535.909
$ws.Range("K132").Value = 6211.399800000001
$ws.Range("L132").Value = 19607.727
$ws.Range("M132").Value = -3681.399800000001
$ws.Range("N132").Value = -24667.727

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
$ws.Range("H16").Value = 66040.125
$ws.Range("I16").Value = 70376.2
$ws.Range("K16").Value = 70376.2
$ws.Range("M16").Value = -70206.2
$ws.Range("H30").Value = 6500
$ws.Range("I30").Value = 6500
$ws.Range("K30").Value = 6500
$ws.Range("M30").Value = -6392
$ws.Range("H82").Value = 1994.0555
$ws.Range("I82").Value = 1383.0834
$ws.Range("K82").Value = 1383.0834
$ws.Range("M82").Value = -1022.0834
$ws.Range("H85").Value = 1994.0555
$ws.Range("I85").Value = 1383.0834
$ws.Range("K85").Value = 1383.0834
$ws.Range("M85").Value = -135.0834

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 1077
$ws.Range("I12").Value = 1077
$ws.Range("K12").Value = 1077
$ws.Range("M12").Value = -935
$ws.Range("H20").Value = 6425
$ws.Range("J20").Value = 12500
$ws.Range("L20").Value = 12500
$ws.Range("N20").Value = -12980
$ws.Range("H29").Value = 5028.4287
$ws.Range("J29").Value = 1200
$ws.Range("L29").Value = 1200
$ws.Range("N29").Value = -1780
$ws.Range("H113").Value = 515.6
$ws.Range("I113").Value = 461.6
$ws.Range("K113").Value = 1384.8
$ws.Range("M113").Value = 785.1999999999998
$ws.Range("H136").Value = 824.75
$ws.Range("I136").Value = 824.75
$ws.Range("K136").Value = 2474.25
$ws.Range("M136").Value = 75.75
